$wb = $excel.ActiveWorkbook

# --- Sheet 1: test_suite ---
$wsSuite = $wb.Worksheets.Item("test_suite")
$wsSuite.Cells.Item(2,1).Value = "TC1BankManagerLogin"
$wsSuite.Cells.Item(3,1).Value = "TC2AddCustomer"
$wsSuite.Cells.Item(4,1).Value = "TC3OpenAccount"
$wsSuite.Columns.Item(1).AutoFit()

# --- Sheet 2: "1" ---
$ws1 = $wb.Worksheets.Item("1")
$ws1.Cells.Item(2,1).Value = "tester1"
$ws1.Cells.Item(2,2).Value = "test2"
$ws1.Cells.Item(2,3).Value = "asas546"

# --- Sheet 3: TC2addCustomer ---
$wsTC2 = $wb.Worksheets.Item("TC2addCustomer")
$wsTC2.Cells.Item(2,1).Value = "anish"
$wsTC2.Cells.Item(2,2).Value = "Sharma"
$wsTC2.Cells.Item(2,3).Value = "dfdsf"
$wsTC2.Cells.Item(2,5).Value = "Y"

$wsTC2.Cells.Item(3,1).Value = "Anita"
$wsTC2.Cells.Item(3,2).Value = "Jacob"
$wsTC2.Cells.Item(3,3).Value = "fgfg"
$wsTC2.Cells.Item(3,5).Value = "N"

$wsTC2.Cells.Item(4,1).Value = "Nita"
$wsTC2.Cells.Item(4,2).Value = "Kulkarni"
$wsTC2.Cells.Item(4,3).Value = "ghgf"
$wsTC2.Cells.Item(4,5).Value = "Y"

$wsTC2.Cells.Item(5,1).Value = "Rajiv"
$wsTC2.Cells.Item(5,2).Value = "Parv"
$wsTC2.Cells.Item(5,3).Value = "hgg"
$wsTC2.Cells.Item(5,5).Value = "Y"

$wsTC2.Cells.Item(6,1).Value = "Ganga"
$wsTC2.Cells.Item(6,2).Value = "River"
$wsTC2.Cells.Item(6,3).Value = "hjh"
$wsTC2.Cells.Item(6,4).Value = "Customer added successfully"

$wsTC2.Cells.Item(7,1).Value = "Dev"
$wsTC2.Cells.Item(7,2).Value = "Sankar"
$wsTC2.Cells.Item(7,3).Value = "mhjk"
$wsTC2.Cells.Item(7,4).Value = "Customer added successfully"

# --- Sheet 4: OpenAccountTest -> TC3OpenAccount ---
$wsOpen = $wb.Worksheets.Item("OpenAccountTest")
$wsOpen.Name = "TC3OpenAccount"
$wsOpen.Cells.Item(2,1).Value = "Anita Jacob"

# --- Selections to match the saved view state ---
$wsSuite.Activate()
$wsSuite.Range("B6").Select()

$ws1.Activate()
$ws1.Range("B8").Select()

$wsTC2.Activate()
$wsTC2.Range("A3:B3").Select()

$wsOpen.Activate()
$wsOpen.Range("C16").Select()
